# "Error Calculations and Plots"
#
# The underlying data table (IDs in column A, numeric error columns B:F)
# had two whole rows removed ("RM 232" and "SC 92") and a handful of
# individual column-F (and a couple column-C) values swapped between a
# real number and a blank/missing cell, simulating a different pattern
# of missing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "RM 232" row (originally row 26) -----------------------
$ws.Rows.Item(26).Delete()

# After that deletion everything below shifts up by one, so the row that
# used to be "SC 92" (originally row 28) is now row 27. Remove it too.
$ws.Rows.Item(27).Delete()

# --- Toggle a handful of column F (and a couple column C) values -------
# between a real number and blank, to match the new missing-data pattern.

# RM 21 (row 6): F was blank -> now has a value
$ws.Cells.Item(6, 6).Value = 16.43

# RM 38 (row 8): F had a value -> now blank
$ws.Cells.Item(8, 6).Value = $null

# RM 125 (row 19): F was blank -> now has a value
$ws.Cells.Item(19, 6).Value = 17.81

# RM 135 (row 21): F had a value -> now blank
$ws.Cells.Item(21, 6).Value = $null

# RM 140 (row 23): F was blank -> now has a value
$ws.Cells.Item(23, 6).Value = 16.48

# SC 5 (row 26 after the deletions above): C had a value -> now blank
$ws.Cells.Item(26, 3).Value = $null

# SC 101 (row 27 after the deletions above): C was blank -> now has a
# value, and F had a value -> now blank
$ws.Cells.Item(27, 3).Value = 10
$ws.Cells.Item(27, 6).Value = $null

# SC 119 (row 29 after the deletions above): C had a value -> now blank,
# and F was blank -> now has a value
$ws.Cells.Item(29, 3).Value = $null
$ws.Cells.Item(29, 6).Value = 18.06
